$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.742.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.583.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.577.17"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.77%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.199"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.30%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.80"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.162.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "614.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.584.27"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.757.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.43%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.59%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -16.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.13"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.94"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.75%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.40"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "648.42"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.22%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.99%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0481"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.33"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.397.39"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0715"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.97"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.25%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.83"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.09%  "
